$wb = $excel.ActiveWorkbook

# --- Update the informational text on "Hoja1" (A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 5.95 = 24184.52 pesos`n✅ 24184.52 pesos = 5.9 = 951.09 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Update the rate figures on "tasas" sheet ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 168
$wsTasas.Range("O10").Value = 4063
$wsTasas.Range("N12").Value = 4099
$wsTasas.Range("O12").Value = 161.199
